# Oppdatert koblingsfilen som innholder oversikt over gamle og nye fylker
# Inserts two new rows into Ark1 (old Sogn og Fjordane / Møre og Romsdal
# entries each get a duplicate "no spaces" variant row, everything below
# shifts down accordingly).
#
# NOTE: the "Moreogromsdal" row is inserted first (so its shared string is
# registered before "Sognogfjordane"'s), then the "Sognogfjordane" row is
# inserted above it -- this reproduces both the final row layout and the
# shared-string insertion order of the original edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 16 (old row 16 "Sor-Trondelag" etc. shifts down to 17+)
# and fill it with the "Moreogromsdal" duplicate of (old) row 15.
$ws.Rows("16").Insert()
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Møreogromsdal"
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = "Helse Midt-Norge"
$ws.Range("G16").Value = 100024

# Insert a new row 15 (old row 15 "Møre og Romsdal" etc. -- now including the
# row just added above -- shifts down to 16+) and fill it with the
# "Sognogfjordane" duplicate of row 14.
$ws.Rows("15").Insert()
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Sognogfjordane"
$ws.Range("C15").Value = 46
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "Helse Vest"
$ws.Range("G15").Value = 100021
